$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N ("Description" shifts from N to O),
# and populate the new "eIDAS RequesterID" column.
$ws.Range("N1").EntireColumn.Insert()
$ws.Range("N1").Value = "eIDAS RequesterID"
$ws.Range("N2").Value = 1234

# Move the stray empty hyperlink-styled cell from F2 to G2.
$ws.Range("F2").Clear()
$ws.Range("G2").Style = "Hyperlink"

# Widen the columns between "Client name (et)" and "Client shortname (en)"
# to fit their content (closest achievable values given column-width rounding).
$ws.Columns.Item(7).ColumnWidth = 16.8
$ws.Columns.Item(8).ColumnWidth = 13.6
$ws.Columns.Item(9).ColumnWidth = 14.3
$ws.Columns.Item(10).ColumnWidth = 23.3
$ws.Columns.Item(11).ColumnWidth = 29.6

# Update the sheet view selection.
$ws.Range("K2").Select()

$wb.Save()
